$d = $word.ActiveDocument

# Use Find to locate the old text, then assign the new text directly to the
# matched Range. Assigning Range.Text (rather than relying on Find's
# ReplaceWith) keeps straight apostrophes/quotes as-is instead of letting
# AutoCorrect turn them into curly quotes.
function Replace-Text($old, $new) {
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, `
                       $true, 1, $false, "", 0) | Out-Null
    if ($rng.Find.Found) {
        $rng.Text = $new
    } else {
        Write-Output "WARNING: text not found: $old"
    }
}

# Title / name change
Replace-Text "SADJI Thylian" "MARINO Murphy"

# Histoire paragraph
Replace-Text "Thylian SADJI est un homme originaire de Samoa où il y a vécu toute son enfance. Depuis tout petit c'est un grand passionné par la lecture il fit des études pour faire carrière dans le domaine, mais il ne réussit pas et se contenta de faire un simple métier à sa portée. Plus jeune, il faisait déjà des petits boulots à côté des études pour aider sa famille pauvre. Il est follement amoureux de la jeune Xana sa voisine, et serait prêt à tout pour l'épouser." `
              "Murphy MARINO est un homme originaire du Salvador où il y a vécu toute son enfance. Depuis tout petit c'est un grand passionné par la nature il fit des études pour faire carrière dans le domaine, il ne réussit pas et pire encore, il fut pendant un temps SDF avant de se reprendre en main. Plus jeune, il se perdit en forêt, ce qui lui causa un traumatisme profond, depuis il a peur de s'aventurer seul dans des endroits sans aucune présence. Il sait qu'il peut toujours compter sur son meilleur ami Nacime, qui l'aidera peu importe la situation."

# Points de vies / Points de Mana
Replace-Text "Points de vies : 14" "Points de vies : 12"
Replace-Text "                 Points de Mana : 6" "                 Points de Mana : 8"

# Langue(s)
Replace-Text "Arabe, Mirandais" "Ourdou, Hindi des Fidji"

# Compétences
Replace-Text "Acrobatie +25%, Comedie +20%, Armes blanches +10 %" "Medecine +25 %, Armes a feu +10 %, Esquive +5 %"
